# Apply the "add 2022-Q3 data" change to the workbook.
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new worksheet named "2022-Q3" right after "总计".
# ------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($null, $zongji)
$newSheet.Name = "2022-Q3"

# Header row (plain text, never numeric-looking -> safe to assign directly)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Numeric-looking text columns (fund code, scale, position, ratio, value) must
# be forced to Text format first so Excel doesn't coerce them into numbers
# (which would e.g. drop the leading zero from fund codes).
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:G5").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "010041"
$newSheet.Range("C2").Value = "嘉实港股优势混合A"
$newSheet.Range("D2").Value = "40.43"
$newSheet.Range("E2").Value = "89.69"
$newSheet.Range("F2").Value = "4.42"
$newSheet.Range("G2").Value = "1.7870"
$newSheet.Range("H2").Value = 5

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "010042"
$newSheet.Range("C3").Value = "嘉实港股优势混合C"
$newSheet.Range("D3").Value = "5.20"
$newSheet.Range("E3").Value = "89.69"
$newSheet.Range("F3").Value = "4.42"
$newSheet.Range("G3").Value = "0.2298"
$newSheet.Range("H3").Value = 5

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "001942"
$newSheet.Range("C4").Value = "前海开源沪港深汇鑫灵活配置混合A"
$newSheet.Range("D4").Value = "0.17"
$newSheet.Range("E4").Value = "87.24"
$newSheet.Range("F4").Value = "4.66"
$newSheet.Range("G4").Value = "0.0079"
$newSheet.Range("H4").Value = 6

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "001943"
$newSheet.Range("C5").Value = "前海开源沪港深汇鑫灵活配置混合C"
$newSheet.Range("D5").Value = "0.09"
$newSheet.Range("E5").Value = "87.24"
$newSheet.Range("F5").Value = "4.66"
$newSheet.Range("G5").Value = "0.0042"
$newSheet.Range("H5").Value = 6

# The Text NumberFormat applied above leaves a visible "@" style on the
# cells; clean it back up to the default General style (a blank cell is
# always General) while the stored values remain text, exactly as in the
# source sheets.
$newSheet.Range("Z100").Copy()
$newSheet.Range("B2:B5").PasteSpecial(-4122)
$newSheet.Range("D2:G5").PasteSpecial(-4122)

# Copy header / first-column formatting from the equivalent "2022-Q2" sheet so
# the new sheet matches the look of its siblings.
$src = $wb.Worksheets.Item("2022-Q2")
$src.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$src.Range("A2:A3").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# PasteSpecial(Formats) only touches formatting, restore the serial values
# it may have overwritten via the repeated copy source.
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3

$newSheet.Range("A1").Select()

# ------------------------------------------------------------------
# 2) Update "总计" with the new 2022-Q3 summary row.
# ------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q3"
$zongji.Range("C2").Value = 4
$zongji.Range("D2").Value = 2.03

# Restore the plain (un-styled) formatting that row 3 has, which the
# row-insert operation did not carry over correctly to the new row 2.
$zongji.Range("A3:D3").Copy()
$zongji.Range("A2:D2").PasteSpecial(-4122)

# Renumber column A (serial index) for the rows that shifted down.
$zongji.Range("A3").Value = 1
$zongji.Range("A4").Value = 2
$zongji.Range("A5").Value = 3
$zongji.Range("A6").Value = 4
$zongji.Range("A7").Value = 5

$zongji.Range("A1").Select()
